$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3858091235160828
$ws.Range("B1").Value = 0.6176436543464661
$ws.Range("D1").Value = 1.407668113708496
$ws.Range("E1").Value = 0.8605092167854309
